$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$itSheet = $wb.Worksheets.Item("IT")

# --- Update the IT sheet's initial-time year value from 2018 to 2020 ---
$itSheet.Range("B2").Value = 2020

# --- Reset the "About" sheet's stored selection back to the top-left cell,
#     since it is no longer the active tab and should not keep its old
#     mid-sheet selection (A21). ---
$aboutSheet.Select()
$aboutSheet.Range("A1").Select()

# --- Make the "IT" sheet the active/selected sheet (was "About") ---
$itSheet.Select()
$itSheet.Range("B3").Select()

$wb.Save()
